$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '27.146.89'
$ws.Range('E2').Value = '  +1.13%  '

# Row 3
$ws.Range('D3').Value = '1.846.69'
$ws.Range('E3').Value = '  +1.79%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.34%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.64%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.005'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4633'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.62%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3703'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.82%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07374'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.25%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8836'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.36%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07913'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.12%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '19.92'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.09%  '

# Row 13
$ws.Range('D13').Value = '1.850.28'
$ws.Range('E13').Value = '  +1.66%  '

# Row 14
$ws.Range('E14').Value = '  +0.98%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.601'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.91%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.89'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.13%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.22%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008937'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.14%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.005'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.38%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.31%  '

# Row 21
$ws.Range('D21').Value = '27.189.29'
$ws.Range('E21').Value = '  +0.69%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.140'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.05%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.75%  '

# Row 24
$ws.Range('D24').Value = '2.128.46'
$ws.Range('E24').Value = '  +2.88%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.91%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.871'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.27%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.44%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.067'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.11%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.24%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.128'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.08%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08886'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.35%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.972'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.39%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7416'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.07%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.468'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.78%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.141'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.09%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.557'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.27%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.081'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.85%  '

# Row 38
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05269'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.87%  '

# Row 39
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01950'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.84%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.977'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.57%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.097'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.35%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5173'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.77%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1636'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.21%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.243'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.01%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4858'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.68%  '

# Row 46
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.005'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.27%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.20'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.74%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '102.92'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.09%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.629'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.05%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06232'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.57%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.60'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.38%  '
